{"js": "// The paragraph \"\u0412\u043d\u0438\u043c\u0430\u043d\u0438\u0435! \u0414\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0444\u0430\u0439\u043b\u043e\u0432\u043e\u0439 \u0441\u0438\u0441\u0442\u0435\u043c\u043e\u0439 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c\n// \u0442\u043e\u043b\u044c\u043a\u043e OS API.\" had a stray \"_GoBack\" bookmark splitting the run\n// \"...\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \" from \"\u0442\u043e\u043b\u044c\u043a\u043e \" into two separate <w:r> elements.\n// Remove the bookmark and re-merge the text back into a single run.\n\n// 1) Drop the leftover _GoBack bookmark (Word inserts this automatically\n//    at the last edit location; it has no visible effect but was left\n//    splitting the two runs).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Re-write the now-adjacent text \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \" as one run so\n//    the two <w:r> runs collapse back into a single run, matching the\n//    original authored text.\nconst target = context.document.body.search(\"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The paragraph \"\u0412\u043d\u0438\u043c\u0430\u043d\u0438\u0435! \u0414\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0444\u0430\u0439\u043b\u043e\u0432\u043e\u0439 \u0441\u0438\u0441\u0442\u0435\u043c\u043e\u0439 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c\n# \u0442\u043e\u043b\u044c\u043a\u043e OS API.\" had a stray \"_GoBack\" bookmark splitting the run\n# \"...\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \" from \"\u0442\u043e\u043b\u044c\u043a\u043e \" into two separate runs.\n# Remove the bookmark and re-merge the text back into a single run.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the leftover _GoBack bookmark (Word drops this automatically at\n#    the last edit location; it has no visible effect but was left\n#    splitting the two runs).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Re-write the now-adjacent text \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \" as one run so\n#    the two runs collapse back into a single run, matching the original\n#    authored text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \"\n$find.Replacement.Text = \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u044c \u0442\u043e\u043b\u044c\u043a\u043e \"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
